$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update values, keep existing styles
$ws.Range("A2").Value = 2310429
$ws.Range("B2").Value = 11185
$ws.Range("C2").Value = 7
$ws.Range("E2").Value = 3.5

# Row 3: update values; C3/D3/E3 also change style from 15 to 14
$ws.Range("A3").Value = 2316494
$ws.Range("B3").Value = 30605
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3.5
$ws.Range("C2").Copy()
$ws.Range("C3:E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 4-10: clear the data, keep formatting intact
$ws.Range("A4:E10").ClearContents()

# Column A cells in rows 4-10 had no explicit style; re-apply "Normal" so the
# empty cell tags are retained in the saved XML instead of being dropped.
$ws.Range("A4:A10").Style = "Normal"

# Update the active selection to E3
$ws.Range("E3").Select()
